$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) figures on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1430
$ws1.Range("F8").Value = 205

# Update the same figures on the "全部类型" sheet (mirrors "展览")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1430
$ws4.Range("F8").Value = 205
